$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E) to (B:F).
# Excel's native column-insert preserves per-cell formatting/style ids for the
# shifted cells and leaves the brand-new column completely blank.
$ws.Range("A1").EntireColumn.Insert()

# The new column B (old column A) holds the segment names; give it a header.
$ws.Range("B1").Value = "segments"

# Match the header's look (bold font / border / centered-top alignment) by
# copying the existing header formatting from C1 onto B1 - this reuses the
# same style id instead of registering a near-duplicate style.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# The shifted segment-name cells (B2:B20) inherited column A's old header-ish
# style; the new layout keeps them unstyled (plain data cells), so strip the
# per-cell formatting back to the workbook default.
$ws.Range("B2:B20").ClearFormats()

# Fill new column A with the numeric segment index (0-based), matching row order
for ($i = 0; $i -le 18; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $i
}

# Apply the same style used by the header cells (bold/border/align) to the new
# index column A2:A20, matching the original column A's appearance.
$ws.Range("C1").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)

$excel.CutCopyMode = 0
